# Applies the "RemoveGroup" -> "DeleteGroup" renaming edit described in the
# commit diff to the single slide of the presentation.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Helper: replace the first occurrence of $find inside the shape's TextRange
# with $replace, using Characters() so only the targeted run's text is
# touched and all existing run/paragraph formatting is preserved.
# (Positional parameters only -- this host's PowerShell subset does not bind
# named `-param value` arguments reliably.)
function Replace-ShapeText($shape, $find, $replace) {
    $tr = $shape.TextFrame.TextRange
    $full = $tr.Text
    $idx = $full.IndexOf($find)
    if ($idx -ge 0) {
        $sub = $tr.Characters($idx + 1, $find.Length)
        $sub.Text = $replace
    }
}

# 1) Shape 8 ("Rectangle 62"): "r:RemoveGroupCommand" becomes "d:DeleteGroup"
#    on its own line, with a new second paragraph reading "Command".
$shape8 = $s.Shapes.Item(8)
$tr8 = $shape8.TextFrame.TextRange
$tr8.Text = "d:DeleteGroup"
$null = $tr8.InsertAfter([char]13 + "Command")

# 2) Shape 25 ("TextBox 27"): "removeGroup()" -> "deleteGroup()"
Replace-ShapeText $s.Shapes.Item(25) "removeGroup" "deleteGroup"

# 3) Shape 26 ("TextBox 28"): ...("removegroup g/Production") -> ("deletegroup g/Production")
Replace-ShapeText $s.Shapes.Item(26) "removegroup" "deletegroup"

# 4) Shape 33 ("Rectangle 62"): ":RemoveGroupCommand" / "Parser" -> ":DeleteGroupCommand" / "Parser"
Replace-ShapeText $s.Shapes.Item(33) "RemoveGroupCommand" "DeleteGroupCommand"

# 5) Shape 43 ("TextBox 45"): parse ("removegroup g/Production") -> parse ("deletegroup g/Production")
Replace-ShapeText $s.Shapes.Item(43) "removegroup" "deletegroup"
